$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.999999451631079
$ws.Range("E2").Value = 0.999999451631079

# Row 3
$ws.Range("D3").Value = 0.5298892219228659
$ws.Range("E3").Value = 0.5298892219228659

# Row 4
$ws.Range("D4").Value = 0.1640915102849946
$ws.Range("E4").Value = 0.1640915102849946

# Row 5
$ws.Range("D5").Value = 0.009908523250282476
$ws.Range("E5").Value = 0.009908523250282476

# Row 6
$ws.Range("D6").Value = 0.9816247461793194
$ws.Range("E6").Value = 0.9816247461793194

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.00005112407847876086
$ws.Range("E7").Value = 0.9999488759215213

# Row 8
$ws.Range("D8").Value = 0.9697894700808201
$ws.Range("E8").Value = 0.03021052991917994

# Row 9
$ws.Range("D9").Value = 0.9852428483625978
$ws.Range("E9").Value = 0.01475715163740221

# Row 10
$ws.Range("D10").Value = 0.9999999999904157
$ws.Range("E10").Value = 0.000000000009584333326984051

# Row 11
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 2.928384780883789
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.9999999999926235
$ws.Range("E12").Value = 0.9999999999926235

# Row 13
$ws.Range("D13").Value = 0.9999999996385267
$ws.Range("E13").Value = 0.9999999996385267

# Row 14
$ws.Range("D14").Value = 0.4308942062822707
$ws.Range("E14").Value = 0.4308942062822707

# Row 15
$ws.Range("D15").Value = 0.0000000000000005592882440660861
$ws.Range("E15").Value = 0.0000000000000005592882440660861

# Row 16
$ws.Range("D16").Value = 0.994915382659006
$ws.Range("E16").Value = 0.994915382659006

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.00000015471022545037
$ws.Range("E17").Value = 0.9999998452897746

# Row 18
$ws.Range("D18").Value = 0.9992763503012601
$ws.Range("E18").Value = 0.0007236496987399477

# Row 19
$ws.Range("D19").Value = 0.9995453213406974
$ws.Range("E19").Value = 0.0004546786593025987

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 6.890167236328125
$ws.Range("G21").Value = 0.6
